$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 6809750
$ws.Range("I132").Value = 8134905.5
$ws.Range("J132").Value = 18326.5
$ws.Range("K132").Value = 24404716.5
$ws.Range("L132").Value = 54979.5
$ws.Range("M132").Value = -24402186.5
$ws.Range("N132").Value = -60039.5

$ws.Range("H135").Value = 35715236
$ws.Range("I135").Value = 651.75
$ws.Range("J135").Value = 125001700
$ws.Range("K135").Value = 5865.75
$ws.Range("L135").Value = 1125015300
$ws.Range("M135").Value = -3330.75
$ws.Range("N135").Value = -1125020370

$ws.Range("H137").Value = 1952.7885
$ws.Range("I137").Value = 1626.1786
$ws.Range("J137").Value = 2333.8333
$ws.Range("K137").Value = 4878.5358
$ws.Range("L137").Value = 7001.499899999999
$ws.Range("M137").Value = -2328.5358
$ws.Range("N137").Value = -12101.4999

$ws.Range("H138").Value = 2096.899
$ws.Range("I138").Value = 947.375
$ws.Range("J138").Value = 2318.494
$ws.Range("K138").Value = 2842.125
$ws.Range("L138").Value = 6955.482
$ws.Range("M138").Value = 2297.875
$ws.Range("N138").Value = -17235.482

$ws.Range("H141").Value = 1173.625
$ws.Range("I141").Value = 1080.6666
$ws.Range("J141").Value = 1452.5
$ws.Range("K141").Value = 3241.9998
$ws.Range("L141").Value = 4357.5
$ws.Range("M141").Value = 1938.0002
$ws.Range("N141").Value = -14717.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9835.214
$ws.Range("I32").Value = 7254.7124
$ws.Range("K32").Value = 7254.7124
$ws.Range("M32").Value = -6967.7124

$ws.Range("H45").Value = 1237.875
$ws.Range("I45").Value = 1117.1666
$ws.Range("J45").Value = 1600
$ws.Range("K45").Value = 1117.1666
$ws.Range("L45").Value = 1600
$ws.Range("M45").Value = -740.1666
$ws.Range("N45").Value = -2354

$ws.Range("H61").Value = 83334830
$ws.Range("I61").Value = 125001020
$ws.Range("J61").Value = 2449.5
$ws.Range("K61").Value = 125001020
$ws.Range("L61").Value = 2449.5
$ws.Range("M61").Value = -125000808
$ws.Range("N61").Value = -2873.5

$ws.Range("H63").Value = 41669010
$ws.Range("I63").Value = 2407.5
$ws.Range("J63").Value = 250002020
$ws.Range("K63").Value = 2407.5
$ws.Range("L63").Value = 250002020
$ws.Range("M63").Value = -1721.5
$ws.Range("N63").Value = -250003392

$ws.Range("H66").Value = 41669010
$ws.Range("I66").Value = 2407.5
$ws.Range("J66").Value = 250002020
$ws.Range("K66").Value = 12037.5
$ws.Range("L66").Value = 1250010100
$ws.Range("M66").Value = -8605.5
$ws.Range("N66").Value = -1250016964

$ws.Range("H74").Value = 1827.0714
$ws.Range("I74").Value = 1231.5834
$ws.Range("K74").Value = 1231.5834
$ws.Range("M74").Value = -357.5834

$ws.Range("H77").Value = 1827.0714
$ws.Range("I77").Value = 1231.5834
$ws.Range("K77").Value = 6157.916999999999
$ws.Range("M77").Value = -1789.916999999999

$ws.Range("H122").Value = 5004.6665
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 5004.6665
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 15013.9995
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -19913.9995

$ws.Range("H136").Value = 83334830
$ws.Range("I136").Value = 125001020
$ws.Range("J136").Value = 2449.5
$ws.Range("K136").Value = 375003060
$ws.Range("L136").Value = 7348.5
$ws.Range("M136").Value = -375000510
$ws.Range("N136").Value = -12448.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 333367680
$ws.Range("I105").Value = 333367680
$ws.Range("K105").Value = 333367680
$ws.Range("M105").Value = -333365933

$ws.Range("H107").Value = 1346.5454
$ws.Range("I107").Value = 1122.1111
$ws.Range("K107").Value = 1122.1111
$ws.Range("M107").Value = 797.8888999999999

$ws.Range("H134").Value = 3884.0303
$ws.Range("I134").Value = 844.3226
$ws.Range("K134").Value = 2532.9678
$ws.Range("M134").Value = 2.032200000000103

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 66667924
$ws.Range("I16").Value = 90910410
$ws.Range("K16").Value = 90910410
$ws.Range("M16").Value = -90910123

$ws.Range("H58").Value = 7775.95
$ws.Range("I58").Value = 1921.6666
$ws.Range("J58").Value = 10284.929
$ws.Range("K58").Value = 1921.6666
$ws.Range("L58").Value = 10284.929
$ws.Range("M58").Value = -1718.6666
$ws.Range("N58").Value = -10690.929

$ws.Range("H99").Value = 2063.923
$ws.Range("I99").Value = 1876
$ws.Range("J99").Value = 2283.1667
$ws.Range("K99").Value = 1876
$ws.Range("L99").Value = 2283.1667
$ws.Range("M99").Value = -378
$ws.Range("N99").Value = -5279.1667

$ws.Range("H105").Value = 734.2222
$ws.Range("I105").Value = 701
$ws.Range("K105").Value = 701
$ws.Range("M105").Value = 1046

$ws.Range("H107").Value = 985.46155
$ws.Range("J107").Value = 2300
$ws.Range("L107").Value = 2300
$ws.Range("N107").Value = -6140

$ws.Range("H113").Value = 66667924
$ws.Range("I113").Value = 90910410
$ws.Range("K113").Value = 90910410
$ws.Range("M113").Value = -90908240

$ws.Range("H126").Value = 2063.923
$ws.Range("I126").Value = 1876
$ws.Range("J126").Value = 2283.1667
$ws.Range("K126").Value = 5628
$ws.Range("L126").Value = 6849.500100000001
$ws.Range("M126").Value = -3158
$ws.Range("N126").Value = -11789.5001

$ws.Range("H136").Value = 7775.95
$ws.Range("I136").Value = 1921.6666
$ws.Range("J136").Value = 10284.929
$ws.Range("K136").Value = 5764.9998
$ws.Range("L136").Value = 30854.787
$ws.Range("M136").Value = -3214.9998
$ws.Range("N136").Value = -35954.787

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 519.74286
$ws.Range("J5").Value = 1223.625
$ws.Range("L5").Value = 3670.875
$ws.Range("N5").Value = -3894.875

$ws.Range("H68").Value = 645.1111
$ws.Range("J68").Value = 547.9091
$ws.Range("L68").Value = 1643.7273
$ws.Range("N68").Value = -3265.7273

$ws.Range("H71").Value = 645.1111
$ws.Range("J71").Value = 547.9091
$ws.Range("L71").Value = 4931.1819
$ws.Range("N71").Value = -13043.1819

$ws.Range("H122").Value = 1061.9143
$ws.Range("J122").Value = 1259.5416
$ws.Range("L122").Value = 11335.8744
$ws.Range("N122").Value = -16235.8744

$ws.Range("H135").Value = 519.74286
$ws.Range("J135").Value = 1223.625
$ws.Range("L135").Value = 11012.625
$ws.Range("N135").Value = -16082.625

$ws.Range("H138").Value = 3346.375
$ws.Range("J138").Value = 3012.9333
$ws.Range("L138").Value = 9038.7999
$ws.Range("N138").Value = -19318.7999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 50003524
$ws.Range("I70").Value = 41670350
$ws.Range("J70").Value = 66669868
$ws.Range("K70").Value = 41670350
$ws.Range("L70").Value = 66669868
$ws.Range("M70").Value = -41670080
$ws.Range("N70").Value = -66670408

$ws.Range("H73").Value = 50003524
$ws.Range("I73").Value = 41670350
$ws.Range("J73").Value = 66669868
$ws.Range("K73").Value = 41670350
$ws.Range("L73").Value = 66669868
$ws.Range("M73").Value = -41669414
$ws.Range("N73").Value = -66671740

$ws.Range("H80").Value = 5016.6665
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 5016.6665
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws.Range("H113").Value = 1325
$ws.Range("I113").Value = 1980
$ws.Range("J113").Value = 997.5
$ws.Range("K113").Value = 1980
$ws.Range("L113").Value = 997.5
$ws.Range("M113").Value = 190
$ws.Range("N113").Value = -5337.5

$ws.Range("H126").Value = 2178.0908
$ws.Range("I126").Value = 1765.7142
$ws.Range("K126").Value = 5297.142599999999
$ws.Range("M126").Value = -2827.142599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 751.2308
$ws.Range("I22").Value = 657.6
$ws.Range("J22").Value = 1063.3334
$ws.Range("K22").Value = 657.6
$ws.Range("L22").Value = 1063.3334
$ws.Range("M22").Value = -362.6
$ws.Range("N22").Value = -1653.3334

$ws.Range("H27").Value = 751.2308
$ws.Range("I27").Value = 657.6
$ws.Range("J27").Value = 1063.3334
$ws.Range("K27").Value = 657.6
$ws.Range("L27").Value = 1063.3334
$ws.Range("M27").Value = -550.6
$ws.Range("N27").Value = -1277.3334

$ws.Range("H68").Value = 1201.2222
$ws.Range("I68").Value = 1172.5714
$ws.Range("K68").Value = 1172.5714
$ws.Range("M68").Value = -423.5714

$ws.Range("H71").Value = 1201.2222
$ws.Range("I71").Value = 1172.5714
$ws.Range("K71").Value = 5862.857
$ws.Range("M71").Value = -2118.857

$ws.Range("H132").Value = 3427.5
$ws.Range("I132").Value = 5799.5
$ws.Range("K132").Value = 17398.5
$ws.Range("M132").Value = -14868.5

$ws.Range("H136").Value = 1968.5714
$ws.Range("I136").Value = 1266.6666
$ws.Range("K136").Value = 3799.9998
$ws.Range("M136").Value = -1249.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 599
$ws.Range("I107").Value = 532
$ws.Range("K107").Value = 1596
$ws.Range("M107").Value = 324

$ws.Range("H122").Value = 9617819
$ws.Range("I122").Value = 13160430
$ws.Range("J122").Value = 2159.2856
$ws.Range("K122").Value = 39481290
$ws.Range("L122").Value = 6477.8568
$ws.Range("M122").Value = -39478840
$ws.Range("N122").Value = -11377.8568

$ws.Range("H132").Value = 2325.1082
$ws.Range("I132").Value = 2000.9688
$ws.Range("K132").Value = 6002.9064
$ws.Range("M132").Value = -3472.9064

$ws.Range("H136").Value = 1510.4348
$ws.Range("J136").Value = 1818.625
$ws.Range("L136").Value = 5455.875
$ws.Range("N136").Value = -10555.875
